# Apply scheduled profit recalculation updates across all leve-profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 3109.125
$ws.Range("J32").Value = 3267.5715
$ws.Range("L32").Value = 3267.5715
$ws.Range("N32").Value = -3919.5715
$ws.Range("H39").Value = 2792.25
$ws.Range("I39").Value = 68
$ws.Range("J39").Value = 7332.6665
$ws.Range("K39").Value = 204
$ws.Range("L39").Value = 21997.9995
$ws.Range("M39").Value = 92
$ws.Range("N39").Value = -22589.9995
$ws.Range("H41").Value = 1569.3077
$ws.Range("I41").Value = 662
$ws.Range("K41").Value = 662
$ws.Range("M41").Value = -222
$ws.Range("H64").Value = 5013.3
$ws.Range("I64").Value = 4681.4443
$ws.Range("J64").Value = 8000
$ws.Range("K64").Value = 4681.4443
$ws.Range("L64").Value = 8000
$ws.Range("M64").Value = -4433.4443
$ws.Range("N64").Value = -8496
$ws.Range("H67").Value = 5013.3
$ws.Range("I67").Value = 4681.4443
$ws.Range("J67").Value = 8000
$ws.Range("K67").Value = 4681.4443
$ws.Range("L67").Value = 8000
$ws.Range("M67").Value = -3823.4443
$ws.Range("N67").Value = -9716
$ws.Range("H98").Value = 1308.6333
$ws.Range("I98").Value = 1284.7931
$ws.Range("K98").Value = 1284.7931
$ws.Range("M98").Value = 213.2068999999999
$ws.Range("H122").Value = 1308.6333
$ws.Range("I122").Value = 1284.7931
$ws.Range("K122").Value = 3854.379300000001
$ws.Range("M122").Value = -1404.379300000001
$ws.Range("H132").Value = 3687.2856
$ws.Range("I132").Value = 3815.1924
$ws.Range("J132").Value = 2024.5
$ws.Range("K132").Value = 11445.5772
$ws.Range("L132").Value = 6073.5
$ws.Range("M132").Value = -8915.5772
$ws.Range("N132").Value = -11133.5
$ws.Range("H135").Value = 5278.0356
$ws.Range("I135").Value = 6432.35
$ws.Range("K135").Value = 57891.15
$ws.Range("M135").Value = -55356.15
$ws.Range("H138").Value = 5460.6704
$ws.Range("J138").Value = 4760.1885
$ws.Range("L138").Value = 14280.5655
$ws.Range("N138").Value = -24560.5655

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2338.5
$ws.Range("I2").Value = 2177
$ws.Range("K2").Value = 2177
$ws.Range("M2").Value = -2064
$ws.Range("H32").Value = 10874.221
$ws.Range("I32").Value = 7599.391
$ws.Range("K32").Value = 7599.391
$ws.Range("M32").Value = -7312.391
$ws.Range("H45").Value = 3159
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 1000
$ws.Range("M45").Value = -623
$ws.Range("H74").Value = 2101.7307
$ws.Range("I74").Value = 1617.8096
$ws.Range("K74").Value = 1617.8096
$ws.Range("M74").Value = -743.8096
$ws.Range("H77").Value = 2101.7307
$ws.Range("I77").Value = 1617.8096
$ws.Range("K77").Value = 8089.048000000001
$ws.Range("M77").Value = -3721.048000000001
$ws.Range("H110").Value = 978.2857
$ws.Range("I110").Value = 804.9474
$ws.Range("K110").Value = 804.9474
$ws.Range("M110").Value = 1240.0526
$ws.Range("H116").Value = 2338.5
$ws.Range("I116").Value = 2177
$ws.Range("K116").Value = 2177
$ws.Range("M116").Value = 117

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2338.5
$ws.Range("I3").Value = 2177
$ws.Range("K3").Value = 2177
$ws.Range("M3").Value = -2063
$ws.Range("H20").Value = 1941.8064
$ws.Range("I20").Value = 1969.6086
$ws.Range("K20").Value = 1969.6086
$ws.Range("M20").Value = -1722.6086
$ws.Range("H99").Value = 3547.5
$ws.Range("I99").Value = 3547.5
$ws.Range("K99").Value = 3547.5
$ws.Range("M99").Value = -2049.5
$ws.Range("H134").Value = 5491.4126
$ws.Range("I134").Value = 6090.886
$ws.Range("K134").Value = 18272.658
$ws.Range("M134").Value = -15737.658

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 405.52173
$ws.Range("J7").Value = 389.9
$ws.Range("L7").Value = 389.9
$ws.Range("N7").Value = -615.9
$ws.Range("H31").Value = 32098.314
$ws.Range("I31").Value = 3484.9285
$ws.Range("K31").Value = 3484.9285
$ws.Range("M31").Value = -3189.9285
$ws.Range("H34").Value = 32098.314
$ws.Range("I34").Value = 3484.9285
$ws.Range("K34").Value = 3484.9285
$ws.Range("M34").Value = -3282.9285
$ws.Range("H58").Value = 3081
$ws.Range("I58").Value = 2726.762
$ws.Range("K58").Value = 2726.762
$ws.Range("M58").Value = -2523.762
$ws.Range("H99").Value = 2749.4
$ws.Range("I99").Value = 2749.6667
$ws.Range("J99").Value = 2747
$ws.Range("K99").Value = 2749.6667
$ws.Range("L99").Value = 2747
$ws.Range("M99").Value = -1251.6667
$ws.Range("N99").Value = -5743
$ws.Range("H107").Value = 409.33334
$ws.Range("I107").Value = 337.14285
$ws.Range("K107").Value = 337.14285
$ws.Range("M107").Value = 1582.85715
$ws.Range("H126").Value = 2749.4
$ws.Range("I126").Value = 2749.6667
$ws.Range("J126").Value = 2747
$ws.Range("K126").Value = 8249.000100000001
$ws.Range("L126").Value = 8241
$ws.Range("M126").Value = -5779.000100000001
$ws.Range("N126").Value = -13181
$ws.Range("H132").Value = 3038.257
$ws.Range("I132").Value = 3083.5293
$ws.Range("K132").Value = 9250.5879
$ws.Range("M132").Value = -6720.5879
$ws.Range("H134").Value = 46430.22
$ws.Range("I134").Value = 47677.09
$ws.Range("K134").Value = 143031.27
$ws.Range("M134").Value = -140496.27
$ws.Range("H136").Value = 3081
$ws.Range("I136").Value = 2726.762
$ws.Range("K136").Value = 8180.286
$ws.Range("M136").Value = -5630.286

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 547.5294
$ws.Range("I2").Value = 95.75
$ws.Range("J2").Value = 949.1111
$ws.Range("K2").Value = 95.75
$ws.Range("L2").Value = 949.1111
$ws.Range("M2").Value = 17.25
$ws.Range("N2").Value = -1175.1111

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6499.5
$ws.Range("I7").Value = 6499.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 6499.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -6387.5
$ws.Range("N7").ClearContents()
$ws.Range("H46").Value = 2059.5
$ws.Range("I46").Value = 1699.4
$ws.Range("J46").Value = 2659.6667
$ws.Range("K46").Value = 1699.4
$ws.Range("L46").Value = 2659.6667
$ws.Range("M46").Value = -1511.4
$ws.Range("N46").Value = -3035.6667
$ws.Range("H68").Value = 2286.5
$ws.Range("I68").Value = 2252.842
$ws.Range("J68").Value = 2499.6667
$ws.Range("K68").Value = 2252.842
$ws.Range("L68").Value = 2499.6667
$ws.Range("M68").Value = -1503.842
$ws.Range("N68").Value = -3997.6667
$ws.Range("H71").Value = 2286.5
$ws.Range("I71").Value = 2252.842
$ws.Range("J71").Value = 2499.6667
$ws.Range("K71").Value = 11264.21
$ws.Range("L71").Value = 12498.3335
$ws.Range("M71").Value = -7520.210000000001
$ws.Range("N71").Value = -19986.3335
$ws.Range("H126").Value = 6499.5
$ws.Range("I126").Value = 6499.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 19498.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -17028.5
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 3895.9119
$ws.Range("I136").Value = 3947
$ws.Range("K136").Value = 11841
$ws.Range("M136").Value = -9291

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H131").Value = 90000
$ws.Range("J131").Value = 90000
$ws.Range("L131").Value = 90000
$ws.Range("N131").Value = -100080
$ws.Range("H132").Value = 4231.4863
$ws.Range("I132").Value = 4168.472
$ws.Range("K132").Value = 12505.416
$ws.Range("M132").Value = -9975.415999999999
$ws.Range("H136").Value = 8837.486000000001
$ws.Range("I136").Value = 8979.733
$ws.Range("J136").Value = 8227.857
$ws.Range("K136").Value = 26939.199
$ws.Range("L136").Value = 24683.571
$ws.Range("M136").Value = -24389.199
$ws.Range("N136").Value = -29783.571
